$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2215.6667
$ws.Range("I33").Value = 1777
$ws.Range("J33").Value = 3751
$ws.Range("K33").Value = 1777
$ws.Range("L33").Value = 3751
$ws.Range("M33").Value = -1548
$ws.Range("N33").Value = -4209
$ws.Range("H76").Value = 7813.0527
$ws.Range("I76").Value = 9432
$ws.Range("K76").Value = 9432
$ws.Range("M76").Value = -9117
$ws.Range("H79").Value = 7813.0527
$ws.Range("I79").Value = 9432
$ws.Range("K79").Value = 9432
$ws.Range("M79").Value = -8340
$ws.Range("H100").Value = 49998.547
$ws.Range("I100").Value = 49997
$ws.Range("J100").Value = 49999.125
$ws.Range("K100").Value = 49997
$ws.Range("L100").Value = 49999.125
$ws.Range("M100").Value = -49456
$ws.Range("N100").Value = -51081.125
$ws.Range("H137").Value = 421484.1
$ws.Range("I137").Value = 770569.1
$ws.Range("J137").Value = 8929
$ws.Range("K137").Value = 2311707.3
$ws.Range("L137").Value = 26787
$ws.Range("M137").Value = -2309157.3
$ws.Range("N137").Value = -31887

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 717
$ws.Range("I2").Value = 1000.3333
$ws.Range("J2").Value = 504.5
$ws.Range("K2").Value = 1000.3333
$ws.Range("L2").Value = 504.5
$ws.Range("M2").Value = -887.3333
$ws.Range("N2").Value = -730.5
$ws.Range("H31").Value = 20592.428
$ws.Range("I31").Value = 5191.1665
$ws.Range("K31").Value = 5191.1665
$ws.Range("M31").Value = -4897.1665
$ws.Range("H74").Value = 9624720
$ws.Range("I74").Value = 15627770
$ws.Range("J74").Value = 19841.3
$ws.Range("K74").Value = 15627770
$ws.Range("L74").Value = 19841.3
$ws.Range("M74").Value = -15626896
$ws.Range("N74").Value = -21589.3
$ws.Range("H77").Value = 9624720
$ws.Range("I77").Value = 15627770
$ws.Range("J77").Value = 19841.3
$ws.Range("K77").Value = 78138850
$ws.Range("L77").Value = 99206.5
$ws.Range("M77").Value = -78134482
$ws.Range("N77").Value = -107942.5
$ws.Range("H97").Value = 1091.25
$ws.Range("I97").Value = 921.75
$ws.Range("J97").Value = 2277.75
$ws.Range("K97").Value = 921.75
$ws.Range("L97").Value = 2277.75
$ws.Range("M97").Value = -425.75
$ws.Range("N97").Value = -3269.75
$ws.Range("H102").Value = 10522.579
$ws.Range("I102").Value = 10522.579
$ws.Range("K102").Value = 10522.579
$ws.Range("M102").Value = -8900.579
$ws.Range("H110").Value = 1335.8334
$ws.Range("I110").Value = 1238.2
$ws.Range("K110").Value = 1238.2
$ws.Range("M110").Value = 806.8
$ws.Range("H116").Value = 717
$ws.Range("I116").Value = 1000.3333
$ws.Range("J116").Value = 504.5
$ws.Range("K116").Value = 1000.3333
$ws.Range("L116").Value = 504.5
$ws.Range("M116").Value = 1293.6667
$ws.Range("N116").Value = -5092.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 717
$ws.Range("I3").Value = 1000.3333
$ws.Range("J3").Value = 504.5
$ws.Range("K3").Value = 1000.3333
$ws.Range("L3").Value = 504.5
$ws.Range("M3").Value = -886.3333
$ws.Range("N3").Value = -732.5
$ws.Range("H94").Value = 1531.3
$ws.Range("I94").Value = 421.27274
$ws.Range("J94").Value = 2888
$ws.Range("K94").Value = 421.27274
$ws.Range("L94").Value = 2888
$ws.Range("M94").Value = 29.72726
$ws.Range("N94").Value = -3790
$ws.Range("H105").Value = 2357.4243
$ws.Range("I105").Value = 1558.6154
$ws.Range("K105").Value = 1558.6154
$ws.Range("M105").Value = 188.3846000000001
$ws.Range("H107").Value = 566.625
$ws.Range("I107").Value = 504.7143
$ws.Range("K107").Value = 504.7143
$ws.Range("M107").Value = 1415.2857
$ws.Range("H134").Value = 209917.7
$ws.Range("I134").Value = 1257.2439
$ws.Range("K134").Value = 3771.7317
$ws.Range("M134").Value = -1236.7317

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 3250
$ws.Range("J21").Value = 3250
$ws.Range("L21").Value = 3250
$ws.Range("N21").Value = -3720
$ws.Range("H25").Value = 5851.5
$ws.Range("I25").Value = 1703.3334
$ws.Range("K25").Value = 1703.3334
$ws.Range("M25").Value = -1529.3334
$ws.Range("H31").Value = 565618.5
$ws.Range("I31").Value = 9035.9
$ws.Range("K31").Value = 9035.9
$ws.Range("M31").Value = -8740.9
$ws.Range("H34").Value = 565618.5
$ws.Range("I34").Value = 9035.9
$ws.Range("K34").Value = 9035.9
$ws.Range("M34").Value = -8833.9
$ws.Range("H107").Value = 2994.6667
$ws.Range("J107").Value = 6249.25
$ws.Range("L107").Value = 6249.25
$ws.Range("N107").Value = -10089.25
$ws.Range("H132").Value = 4912.923
$ws.Range("I132").Value = 2509.0435
$ws.Range("K132").Value = 7527.130500000001
$ws.Range("M132").Value = -4997.130500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1075.4117
$ws.Range("J107").Value = 1111.375
$ws.Range("L107").Value = 3334.125
$ws.Range("N107").Value = -7174.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 175.58333
$ws.Range("J2").Value = 134.33333
$ws.Range("L2").Value = 134.33333
$ws.Range("N2").Value = -360.33333
$ws.Range("H126").Value = 3955.3333
$ws.Range("I126").Value = 3541.25
$ws.Range("J126").Value = 4428.5713
$ws.Range("K126").Value = 10623.75
$ws.Range("L126").Value = 13285.7139
$ws.Range("M126").Value = -8153.75
$ws.Range("N126").Value = -18225.7139
$ws.Range("H132").Value = 28573806
$ws.Range("I132").Value = 30305370
$ws.Range("J132").Value = 3007
$ws.Range("K132").Value = 90916110
$ws.Range("L132").Value = 9021
$ws.Range("M132").Value = -90913580
$ws.Range("N132").Value = -14081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3224.9167
$ws.Range("I22").Value = 3112.375
$ws.Range("J22").Value = 3450
$ws.Range("K22").Value = 3112.375
$ws.Range("L22").Value = 3450
$ws.Range("M22").Value = -2817.375
$ws.Range("N22").Value = -4040
$ws.Range("H27").Value = 3224.9167
$ws.Range("I27").Value = 3112.375
$ws.Range("J27").Value = 3450
$ws.Range("K27").Value = 3112.375
$ws.Range("L27").Value = 3450
$ws.Range("M27").Value = -3005.375
$ws.Range("N27").Value = -3664
$ws.Range("H93").Value = 29414030
$ws.Range("I93").Value = 83334810
$ws.Range("J93").Value = 2696.4092
$ws.Range("K93").Value = 83334810
$ws.Range("L93").Value = 2696.4092
$ws.Range("M93").Value = -83333562
$ws.Range("N93").Value = -5192.4092
$ws.Range("H100").Value = 3422.6365
$ws.Range("I100").Value = 2650
$ws.Range("K100").Value = 2650
$ws.Range("M100").Value = -2109
$ws.Range("H132").Value = 392013.66
$ws.Range("I132").Value = 8379.637000000001
$ws.Range("K132").Value = 25138.911
$ws.Range("M132").Value = -22608.911
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H96").Value = 1892
$ws.Range("I96").Value = 1561
$ws.Range("J96").Value = 2146.6155
$ws.Range("K96").Value = 1561
$ws.Range("L96").Value = 2146.6155
$ws.Range("M96").Value = -188
$ws.Range("N96").Value = -4892.6155
$ws.Range("H107").Value = 35715960
$ws.Range("I107").Value = 55557160
$ws.Range("J107").Value = 1796
$ws.Range("K107").Value = 166671480
$ws.Range("L107").Value = 5388
$ws.Range("M107").Value = -166669560
$ws.Range("N107").Value = -9228
$ws.Range("H126").Value = 1953.5
$ws.Range("I126").Value = 1604.6666
$ws.Range("K126").Value = 4813.9998
$ws.Range("M126").Value = -2343.9998
$ws.Range("H132").Value = 316606.62
$ws.Range("I132").Value = 3658.423
$ws.Range("K132").Value = 10975.269
$ws.Range("M132").Value = -8445.269
